# Add in system control and basic PWM control
#
# Fills in the "Mailboxes (CANB)" sheet with Input/Output, CPU and
# Location data for mailboxes 4-7 (rows 7-10), introducing two new
# shared strings ("HO_CAN->modes" and "HO_CAN->motorControlSlow") along
# the way, and leaves the selection on H11 (matching where the author's
# cursor ended up after typing in the new data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mailboxes (CANB)")
$ws.Activate()

# Mailbox 4 (row 7) - Output on CPU1, raw CAN inputs
$ws.Range("C7").Value = "Output"
$ws.Range("D7").Value = "CPU1"
$ws.Range("E7").Value = "HO_CAN->rawInputs"

# Mailbox 5 (row 8) - Output on CPU1, raw CAN inputs
$ws.Range("C8").Value = "Output"
$ws.Range("D8").Value = "CPU1"
$ws.Range("E8").Value = "HO_CAN->rawInputs"

# Mailbox 6 (row 9) - Output on CPU1, system control modes
$ws.Range("C9").Value = "Output"
$ws.Range("D9").Value = "CPU1"
$ws.Range("E9").Value = "HO_CAN->modes"

# Mailbox 7 (row 10) - Output on CPU1, basic PWM / motor control
$ws.Range("C10").Value = "Output"
$ws.Range("D10").Value = "CPU1"
$ws.Range("E10").Value = "HO_CAN->motorControlSlow"

# Leave the selection where the author's cursor ended up
$ws.Range("H11").Select()
